# Weekly update: insert a new price-report row for Membrillo (Terminal La
# Palmera de La Serena) above the existing row 43, pushing the historical
# rows 43-64 down to 44-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 43 (shifts old rows 43:64 -> 44:65).
$ws.Rows.Item(43).Insert()

# Fill the new row 43 with the latest weekly record.
$ws.Cells.Item(43, 1).Value = 8
$ws.Cells.Item(43, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(43, 3).Value = "Coquimbo"
$ws.Cells.Item(43, 4).Value = 45029
$ws.Cells.Item(43, 5).Value = 4
$ws.Cells.Item(43, 6).Value = "Fruta"
$ws.Cells.Item(43, 7).Value = 100104
$ws.Cells.Item(43, 8).Value = "Frutos de pepita"
$ws.Cells.Item(43, 9).Value = 100104003
$ws.Cells.Item(43, 10).Value = "Membrillo"
$ws.Cells.Item(43, 11).Value = "Champion"
$ws.Cells.Item(43, 12).Value = "Primera"
$ws.Cells.Item(43, 13).Value = 16
$ws.Cells.Item(43, 14).Value = 290000
$ws.Cells.Item(43, 15).Value = 300000
$ws.Cells.Item(43, 16).Value = 295000
$ws.Cells.Item(43, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(43, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(43, 19).Value = 656
$ws.Cells.Item(43, 20).Value = 450

# Apply the same date number-format used by the rest of column D so the
# new cell renders as a date rather than a raw serial number.
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
